# ELS_stimbreakdown.xlsx — "updating syncopation stims w/ degrees"
#
# The stim table is re-sorted by syncopation degree (ascending), and the
# meaning of columns C/D is swapped: C now holds the syncopation_group
# label (low/medium/mediun/high) and D now holds the numeric
# syncopation_degree. The header row and shared-string table are updated
# to match (the old "syncopation_degree" header string is dropped and a
# new "syncopation_level" header string is used for D1); "stim_no" (A1)
# and "syncopation_group" (C1) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ----------------------------------------------------------
$ws.Range("A1").Value = "stim_no"
$ws.Range("B1").Value = "stim_name"
$ws.Range("C1").Value = "syncopation_group"
$ws.Range("D1").Value = "syncopation_level"

# --- Data rows (re-sorted ascending by syncopation degree) ---------------
# row, stim_no, stim_name, syncopation_group, syncopation_degree
$data = @(
    @(1,  "Experimenter-composed Low no. 1.mp3",                          "low",    6),
    @(2,  "Experimenter-composed Low no. 6.mp3",                          "low",    8),
    @(3,  "Honeydippers - Impeach the president.mp3",                     "low",    9),
    @(4,  "Annette Peacock - Survival.mp3",                               "low",    12),
    @(5,  "Kool and the Gang - Jungle jazz.mp3",                          "low",    14),
    @(6,  "Gaz - Sing Sing.mp3",                                          "medium", 17),
    @(7,  "The Turtles - I'm chief Kamanawalalea.mp3",                    "medium", 19),
    @(8,  "Please - Sing a simple song.mp3",                              "medium", 29),
    @(9,  "Ike and Tina Turner - Cussin' and cryin' and carryin' on.mp3", "medium", 32),
    @(10, "Lou Donaldson - Ode to Billy Joe.mp3",                         "mediun", 45),
    @(11, "Experimenter-composed High no. 3.mp3",                         "high",   58),
    @(12, "Experimenter-composed High no. 4.mp3",                         "high",   62),
    @(13, "Experimenter-composed High no. 6.mp3",                         "high",   65),
    @(14, "Experimenter-composed High no. 5.mp3",                         "high",   78),
    @(15, "Experimenter-composed High no. 2.mp3",                         "high",   81)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# --- Formatting ------------------------------------------------------------
# Column B widened to fit the longest stim_name (bestFit-style AutoFit width).
$ws.Columns.Item(2).ColumnWidth = 49.5

# Selection moved to the whole of column C.
[void]$ws.Range("C1:C1048576").Select()

# The old manual-sort record no longer matches the (now re-sorted) data, so
# Excel drops it on save - clear the range's sort state to match.
$ws.Sort.SortFields.Clear()
